# "changed u-net dimensions to volumetric"
#
# The workbook is a U-Net receptive-field / output-size calculator. Every
# changed <v> in the diff (B3, R3, T3, B4, T4, D5, P5, F7, N7, H9, L9) is a
# *formula result* that cascades from a single upstream input: the bottom
# "Input" cell J11, which goes from 28 to 40 (2D patch -> volumetric patch).
# Updating J11 and letting Excel recalculate reproduces every one of those
# values automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The single real data edit driving the whole recalculation cascade.
$ws.Range("J11").Value = 40

# Leave the selection where the author finished editing.
$ws.Range("J12").Select()
